$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete row 26 ("RM 232") entirely - remaining rows shift up
$ws.Rows.Item(26).Delete()

# After the first deletion, the row that was "SC 92" (old row 28) is now at row 27.
# Delete it too - remaining rows shift up again
$ws.Rows.Item(27).Delete()

# Fix up the cells whose imputed values changed as part of this edit
$ws.Range("E26").Value = -5
$ws.Range("E27").Value = ""
$ws.Range("F33").Value = 17.53
